$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H18").Value = 14399
$ws.Range("J18").Value = 16248.75
$ws.Range("L18").Value = 16248.75
$ws.Range("N18").Value = -16816.75

$ws.Range("H28").Value = 301.5
$ws.Range("I28").Value = 243.77777
$ws.Range("J28").Value = 405.4
$ws.Range("K28").Value = 243.77777
$ws.Range("L28").Value = 405.4
$ws.Range("M28").Value = 241.22223
$ws.Range("N28").Value = -1375.4

$ws.Range("H76").Value = 4692397.5
$ws.Range("I76").Value = 7814131
$ws.Range("K76").Value = 7814131
$ws.Range("M76").Value = -7813816

$ws.Range("H79").Value = 4692397.5
$ws.Range("I79").Value = 7814131
$ws.Range("K79").Value = 7814131
$ws.Range("M79").Value = -7813039

$ws.Range("H108").Value = 67999
$ws.Range("J108").Value = 67999
$ws.Range("L108").Value = 67999
$ws.Range("N108").Value = -75679

$ws.Range("H129").Value = 870.7091
$ws.Range("J129").Value = 869.78
$ws.Range("L129").Value = 2609.34
$ws.Range("N129").Value = -12609.34

$ws.Range("H131").Value = 2559.2104
$ws.Range("I131").Value = 902.7143
$ws.Range("J131").Value = 3525.5
$ws.Range("K131").Value = 2708.1429
$ws.Range("L131").Value = 10576.5
$ws.Range("M131").Value = 2331.8571
$ws.Range("N131").Value = -20656.5

$ws.Range("H132").Value = 1288.091
$ws.Range("I132").Value = 1130.0667
$ws.Range("J132").Value = 2868.3333
$ws.Range("K132").Value = 3390.2001
$ws.Range("L132").Value = 8604.999899999999
$ws.Range("M132").Value = -860.2001
$ws.Range("N132").Value = -13664.9999

$ws.Range("H135").Value = 671.6111
$ws.Range("I135").Value = 578.1429
$ws.Range("J135").Value = 998.75
$ws.Range("K135").Value = 5203.2861
$ws.Range("L135").Value = 8988.75
$ws.Range("M135").Value = -2668.2861
$ws.Range("N135").Value = -14058.75

$ws.Range("H137").Value = 2039.8667
$ws.Range("J137").Value = 2324.75
$ws.Range("L137").Value = 6974.25
$ws.Range("N137").Value = -12074.25

$ws.Range("H138").Value = 2395.2542
$ws.Range("I138").Value = 2828.6365
$ws.Range("J138").Value = 2137.5676
$ws.Range("K138").Value = 8485.9095
$ws.Range("L138").Value = 6412.702799999999
$ws.Range("M138").Value = -3345.9095
$ws.Range("N138").Value = -16692.7028

$ws.Range("H141").Value = 3097.6206
$ws.Range("I141").Value = 2435.739
$ws.Range("K141").Value = 7307.217000000001
$ws.Range("M141").Value = -2127.217000000001

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H2").Value = 1246.375
$ws.Range("I2").Value = 1253.5
$ws.Range("J2").Value = 1225
$ws.Range("K2").Value = 1253.5
$ws.Range("L2").Value = 1225
$ws.Range("M2").Value = -1140.5
$ws.Range("N2").Value = -1451

$ws.Range("H32").Value = 2361.4087
$ws.Range("I32").Value = 1602.0513
$ws.Range("K32").Value = 1602.0513
$ws.Range("M32").Value = -1315.0513

$ws.Range("H61").Value = 2861.7368
$ws.Range("I61").Value = 1508.5
$ws.Range("J61").Value = 5181.5713
$ws.Range("K61").Value = 1508.5
$ws.Range("L61").Value = 5181.5713
$ws.Range("M61").Value = -1296.5
$ws.Range("N61").Value = -5605.5713

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

$ws.Range("H97").Value = 955.3333
$ws.Range("I97").Value = 799.7143
$ws.Range("K97").Value = 799.7143
$ws.Range("M97").Value = -303.7143

$ws.Range("H102").Value = 2563.625
$ws.Range("I102").Value = 2251.5
$ws.Range("K102").Value = 2251.5
$ws.Range("M102").Value = -629.5

$ws.Range("H110").Value = 1716.3846
$ws.Range("I110").Value = 1025
$ws.Range("K110").Value = 1025
$ws.Range("M110").Value = 1020

$ws.Range("H116").Value = 1246.375
$ws.Range("I116").Value = 1253.5
$ws.Range("J116").Value = 1225
$ws.Range("K116").Value = 1253.5
$ws.Range("L116").Value = 1225
$ws.Range("M116").Value = 1040.5
$ws.Range("N116").Value = -5813

$ws.Range("H132").Value = 1405.8125
$ws.Range("I132").Value = 1133.5869
$ws.Range("K132").Value = 3400.7607
$ws.Range("M132").Value = -870.7606999999998

$ws.Range("H136").Value = 2861.7368
$ws.Range("I136").Value = 1508.5
$ws.Range("J136").Value = 5181.5713
$ws.Range("K136").Value = 4525.5
$ws.Range("L136").Value = 15544.7139
$ws.Range("M136").Value = -1975.5
$ws.Range("N136").Value = -20644.7139

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H3").Value = 1246.375
$ws.Range("I3").Value = 1253.5
$ws.Range("J3").Value = 1225
$ws.Range("K3").Value = 1253.5
$ws.Range("L3").Value = 1225
$ws.Range("M3").Value = -1139.5
$ws.Range("N3").Value = -1453

$ws.Range("H86").Value = 336066.66
$ws.Range("J86").Value = 502600
$ws.Range("L86").Value = 502600
$ws.Range("N86").Value = -504846

$ws.Range("H89").Value = 336066.66
$ws.Range("J89").Value = 502600
$ws.Range("L89").Value = 2513000
$ws.Range("N89").Value = -2524232

$ws.Range("H128").Value = 3166.6667
$ws.Range("I128").Value = 3166.6667
$ws.Range("K128").Value = 9500.000100000001
$ws.Range("M128").Value = -7010.000100000001

$ws.Range("H134").Value = 6990.304
$ws.Range("I134").Value = 8757.471
$ws.Range("J134").Value = 1983.3334
$ws.Range("K134").Value = 26272.413
$ws.Range("L134").Value = 5950.0002
$ws.Range("M134").Value = -23737.413
$ws.Range("N134").Value = -11020.0002

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 1828.2084
$ws.Range("I31").Value = 1483.2307
$ws.Range("J31").Value = 2235.9092
$ws.Range("K31").Value = 1483.2307
$ws.Range("L31").Value = 2235.9092
$ws.Range("M31").Value = -1188.2307
$ws.Range("N31").Value = -2825.9092

$ws.Range("H34").Value = 1828.2084
$ws.Range("I34").Value = 1483.2307
$ws.Range("J34").Value = 2235.9092
$ws.Range("K34").Value = 1483.2307
$ws.Range("L34").Value = 2235.9092
$ws.Range("M34").Value = -1281.2307
$ws.Range("N34").Value = -2639.9092

$ws.Range("H53").Value = 68374.5
$ws.Range("J53").Value = 68374.5
$ws.Range("L53").Value = 68374.5
$ws.Range("N53").Value = -69588.5

$ws.Range("H58").Value = 1662.862
$ws.Range("I58").Value = 1009.7619
$ws.Range("J58").Value = 3377.25
$ws.Range("K58").Value = 1009.7619
$ws.Range("L58").Value = 3377.25
$ws.Range("M58").Value = -806.7619
$ws.Range("N58").Value = -3783.25

$ws.Range("H62").Value = 2592.3333
$ws.Range("J62").Value = 2475
$ws.Range("L62").Value = 2475
$ws.Range("N62").Value = -3723

$ws.Range("H65").Value = 2592.3333
$ws.Range("J65").Value = 2475
$ws.Range("L65").Value = 12375
$ws.Range("N65").Value = -18615

$ws.Range("H132").Value = 2572.423
$ws.Range("I132").Value = 1613.1538
$ws.Range("J132").Value = 3531.6924
$ws.Range("K132").Value = 4839.4614
$ws.Range("L132").Value = 10595.0772
$ws.Range("M132").Value = -2309.4614
$ws.Range("N132").Value = -15655.0772

$ws.Range("H134").Value = 2466.1667
$ws.Range("I134").Value = 2087
$ws.Range("J134").Value = 5499.5
$ws.Range("K134").Value = 6261
$ws.Range("L134").Value = 16498.5
$ws.Range("M134").Value = -3726
$ws.Range("N134").Value = -21568.5

$ws.Range("H136").Value = 1662.862
$ws.Range("I136").Value = 1009.7619
$ws.Range("J136").Value = 3377.25
$ws.Range("K136").Value = 3029.2857
$ws.Range("L136").Value = 10131.75
$ws.Range("M136").Value = -479.2856999999999
$ws.Range("N136").Value = -15231.75

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H5").Value = 713.5
$ws.Range("J5").Value = 741.3333
$ws.Range("L5").Value = 2223.9999
$ws.Range("N5").Value = -2447.9999

$ws.Range("H7").Value = 512.5
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 640
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 1920
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -2144

$ws.Range("H22").Value = 12500.25
$ws.Range("I22").Value = 12500.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 37500.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -37331.75
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 12500.25
$ws.Range("I27").Value = 12500.25
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 37500.75
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -37398.75
$ws.Range("N27").ClearContents()

$ws.Range("H92").Value = 325
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -3396

$ws.Range("H113").Value = 9790.091
$ws.Range("J113").Value = 887.5
$ws.Range("L113").Value = 2662.5
$ws.Range("N113").Value = -7002.5

$ws.Range("H131").Value = 794.33673
$ws.Range("I131").Value = 525.4
$ws.Range("J131").Value = 808.7957
$ws.Range("K131").Value = 1576.2
$ws.Range("L131").Value = 2426.3871
$ws.Range("M131").Value = 3463.8
$ws.Range("N131").Value = -12506.3871

$ws.Range("H135").Value = 713.5
$ws.Range("J135").Value = 741.3333
$ws.Range("L135").Value = 6671.9997
$ws.Range("N135").Value = -11741.9997

$ws.Range("H140").Value = 1820.119
$ws.Range("I140").Value = 1033
$ws.Range("J140").Value = 2355.36
$ws.Range("K140").Value = 3099
$ws.Range("L140").Value = 7066.08
$ws.Range("M140").Value = 2081
$ws.Range("N140").Value = -17426.08

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H70").Value = 4713.857
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 4332.8335
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 4332.8335
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -4872.8335

$ws.Range("H73").Value = 4713.857
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 4332.8335
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 4332.8335
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -6204.8335

$ws.Range("H97").Value = 1794.1
$ws.Range("I97").Value = 1832.8572
$ws.Range("J97").Value = 1703.6666
$ws.Range("K97").Value = 1832.8572
$ws.Range("L97").Value = 1703.6666
$ws.Range("M97").Value = -1336.8572
$ws.Range("N97").Value = -2695.6666

$ws.Range("H102").Value = 3588.6667
$ws.Range("I102").Value = 3883
$ws.Range("K102").Value = 3883
$ws.Range("M102").Value = -2261

$ws.Range("H103").Value = 30302
$ws.Range("J103").Value = 30302
$ws.Range("L103").Value = 30302
$ws.Range("N103").Value = -32646

$ws.Range("H132").Value = 3327
$ws.Range("I132").Value = 3210.8667
$ws.Range("K132").Value = 9632.6001
$ws.Range("M132").Value = -7102.6001

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H46").Value = 2938.2222
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 3130.5
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 3130.5
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -3506.5

$ws.Range("H68").Value = 2157.3333
$ws.Range("I68").Value = 1916.5714
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1916.5714
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1167.5714
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 2157.3333
$ws.Range("I71").Value = 1916.5714
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 9582.857
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -5838.857
$ws.Range("N71").Value = -22488

$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248

$ws.Range("H122").Value = 4833.6875
$ws.Range("I122").Value = 3689.2
$ws.Range("K122").Value = 11067.6
$ws.Range("M122").Value = -8617.599999999999

$ws.Range("H132").Value = 4469.25
$ws.Range("I132").Value = 2964.5
$ws.Range("J132").Value = 5974
$ws.Range("K132").Value = 8893.5
$ws.Range("L132").Value = 17922
$ws.Range("M132").Value = -6363.5
$ws.Range("N132").Value = -22982

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H108").Value = 65998
$ws.Range("J108").Value = 65998
$ws.Range("L108").Value = 65998
$ws.Range("N108").Value = -73678

$ws.Range("H132").Value = 1659
$ws.Range("I132").Value = 1308.6666
$ws.Range("J132").Value = 2499.8
$ws.Range("K132").Value = 3925.9998
$ws.Range("L132").Value = 7499.400000000001
$ws.Range("M132").Value = -1395.9998
$ws.Range("N132").Value = -12559.4

$ws.Range("H136").Value = 2805.9565
$ws.Range("I136").Value = 2322.6
$ws.Range("K136").Value = 6967.799999999999
$ws.Range("M136").Value = -4417.799999999999
